# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 04:05"

# Brasil (row 5) - updated totals
$ws.Range("B5").Value = 394507
$ws.Range("C5").Value = 2147
$ws.Range("E5").Value = 211321
$ws.Range("G5").Value = 44
$ws.Range("H5").Value = 24593

# Honduras (row 71) - updated totals
$ws.Range("B71").Value = 4401
$ws.Range("C71").Value = 212
$ws.Range("D71").Value = 493
$ws.Range("E71").Value = 3720
$ws.Range("G71").Value = 6
$ws.Range("H71").Value = 188

# Guatemala moves above Hungria in the country ranking (rows 74-75 swap
# identities); row 74 becomes Guatemala with refreshed totals, row 75
# becomes Hungria retaining its previous totals.
$ws.Range("A74").Value = "Guatemala"
$ws.Range("B74").Value = 3954
$ws.Range("C74").Value = 194
$ws.Range("D74").Value = 289
$ws.Range("E74").Value = 3602
$ws.Range("G74").Value = 4
$ws.Range("H74").Value = 63

$ws.Range("A75").Value = "Hungria"
$ws.Range("B75").Value = 3771
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 1836
$ws.Range("E75").Value = 1436
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 499

# Nueva Zelanda (row 98)
$ws.Range("D98").Value = 1462
$ws.Range("E98").Value = 21

# Tunez (row 111)
$ws.Range("D111").Value = 929
$ws.Range("E111").Value = 74
